$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New red-bold font style used for the "important" row 14
$ws.Range("B14:D14").Font.Bold = $true
$ws.Range("B14:D14").Font.Color = 255

$ws.Range("E14").Font.Bold = $true
$ws.Range("E14").Font.Color = 255
$ws.Range("E14").WrapText = $true

# Row 14 values
$ws.Range("B14").Value = 24
$ws.Range("C14").Value = "5. Camada de Serviço"
$ws.Range("D14").Value = "24. Gerenciamento de Transações"
$ws.Range("E14").Value = "0:36 - IMPORTANTE: anotação @Transactional; é utilizada tanto sobre a assinatura de uma classe como a assinatura de um método. Para gerenciar transaçoes pelo spring é utilizada esta anotação. A anotação possui um atributo chamada ""readonly""... um booleano que define se vai precisar que uma transação seja aberta ou não. Quando uma transação é aberta, isso bloqueia o a tabela para outros usuários para operações como metodos de insert, update, delete. Portanto isso pode ser interessante para definir em metodos somente de leitura de dados."

$ws.Rows("14").RowHeight = 120

# Row 15 values
$ws.Range("B15").Value = 25
$ws.Range("C15").Value = "5. Camada de Serviço"
$ws.Range("D15").Value = "25. Incluindo Services"
$ws.Range("E15").Value = "nenhuma anotação na aula porém foi o inicio da implementação dos serviçes, das interfaces e como elas agem. Interessante revisar caso sugir duvidas"

$ws.Range("E15").WrapText = $true

$ws.Rows("15").RowHeight = 30

$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D21").Select()
